$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (columns D, J, K, L, M, O, P) after re-shuffling the
# weekly price records for "Zapallo italiano" / Mapocho Venta Directa de Santiago.

$rows = @{
    2  = @{ D = "01/24/2022"; J = 30; K = 11000; L = 11000; M = 11000; O = "Provincia de Limarí"; P = 183 }
    3  = @{ D = "12/21/2020"; J = 15; K = 7000;  L = 7000;  M = 7000;  O = "Provincia de Limarí"; P = 117 }
    4  = @{ D = "04/26/2021"; J = 30; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí"; P = 167 }
    6  = @{ D = "01/04/2021"; J = 10; K = 9000;  L = 9000;  M = 9000;  O = "Provincia de Limarí"; P = 150 }
    7  = @{ D = "05/17/2021"; J = 25; K = 10000; L = 11000; M = 10400; O = "Provincia de Limarí"; P = 173 }
    8  = @{ D = "12/14/2020"; J = 15; K = 7000;  L = 7000;  M = 7000;  O = "Provincia de Limarí"; P = 117 }
    9  = @{ D = "04/29/2021"; J = 25; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí"; P = 167 }
    11 = @{ D = "04/05/2021"; J = 20; K = 9000;  L = 9000;  M = 9000;  O = "Provincia de Limarí"; P = 150 }
    12 = @{ D = "02/16/2021"; J = 80; K = 10000; L = 11000; M = 10375; O = "Provincia de Quillota"; P = 173 }
    13 = @{ D = "03/29/2021"; J = 35; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí"; P = 167 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
}
